$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was empty) + replace the duplicated "Contact" row with a
# "Jurisdiction" row, then drop the extra duplicate "Contact" row entirely.
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"
$meta.Rows.Item(11).Delete()

# --- Elements sheet ---
$elem = $wb.Worksheets.Item("Elements")

# Root Extension row: Short/Definition now mirror the profile Title/Description
$elem.Range("K2").Value = "Health Data Connect PCP Responsibility Indicator"
$elem.Range("L2").Value = "Indicates whether the primary care physician is the physician considered either responsible or accountable for this claim"

# Column K widened to fit the new (longer) Short text
$elem.Columns.Item(11).ColumnWidth = 45.42
